$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename the column headers: "<name>_old" -> "<name>_FV2404" and
#    "<name>_new" -> "<name>_FV2410" (columns A1:U1).
# ---------------------------------------------------------------------------
$headerRange = $ws.Range("A1:U1")
for ($c = 1; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value2
    if ($val -ne $null) {
        $newVal = $val -replace "_old$", "_FV2404"
        $newVal = $newVal -replace "_new$", "_FV2410"
        if ($newVal -ne $val) {
            $cell.Value = $newVal
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Turn the used range A1:U57 into an Excel Table ("Table1") with an
#    AutoFilter, keeping the header row's existing formatting (bold, fill,
#    border) as plain cell formatting rather than a captured header-row
#    differential format. We do this by stashing the header formatting in a
#    scratch range, resetting the header cells to Normal so the table
#    creation has nothing special to capture, then restoring the original
#    formatting once the table exists.
# ---------------------------------------------------------------------------
$scratchRange = $ws.Range("A100:U100")
$headerRange.Copy()
$scratchRange.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$headerRange.Style = "Normal"

$tableRange = $ws.Range("A1:U57")
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

$scratchRange.Copy()
$headerRange.PasteSpecial(-4122) | Out-Null    # xlPasteFormats
$scratchRange.Clear() | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3) Freeze the header row (split below row 1).
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
